# feat(module) : add rules
#
# Turns the per-day "work note" cells in column C of the "this week" table
# (rows 10-14) into templated placeholders {thisWeekdo1}..{thisWeekdo5},
# gives the "next week" table's column C (rows 17-21) five distinct
# placeholders {nextWeekdo1}..{nextWeekdo5} instead of one repeated literal
# string, and parameterises the closing work-summary sentence in A23 with a
# {summer} placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "This week" table: column C was empty; now carries a placeholder per
#     weekday row, formatted as Text like its sibling cells (B5/B6/B7 etc.).
$thisWeekDoCells = @("C10", "C11", "C12", "C13", "C14")
for ($i = 0; $i -lt $thisWeekDoCells.Length; $i++) {
    $cell = $ws.Range($thisWeekDoCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = "{thisWeekdo" + ($i + 1) + "}"
}

# --- "Next week" table: column C used to repeat the same literal string on
#     every row; now each row gets its own placeholder.
$nextWeekDoCells = @("C17", "C18", "C19", "C20", "C21")
for ($i = 0; $i -lt $nextWeekDoCells.Length; $i++) {
    $ws.Range($nextWeekDoCells[$i]).Value = "{nextWeekdo" + ($i + 1) + "}"
}

# --- Closing summary sentence: swap the hard-coded tail for a placeholder.
$ws.Range("A23").Value = "本周工作总结：{summer}"
